{"js": "// Office.js (Word JavaScript API) script.\n// This script is the body of: async (context) => { ... }\n// It applies a series of precise, unique find-and-replace edits\n// to the body text of the document, matching the target diff.\n\nconst body = context.document.body;\n\n// Ordered list of [searchText, replacementText] pairs. Each search string\n// is unique within the document, so we can safely target the first (only)\n// hit returned by body.search().\nconst edits = [\n  // Paragraph: \"Everything we hear or see, ...\"\n  [\n    \"school, or that burst\",\n    \"school or that burst\",\n  ],\n  [\n    \"significantly, and has leveraged our dependence of technology, completely ignoring the down-sides\",\n    \"significantly and has leveraged our dependence on technology, completely ignoring the downsides\",\n  ],\n  [\n    \"All of this, increases the likelihood of transforming media into ones liking and thereby causing psychological manipulation of the society as a whole.\",\n    \"All of this increases the likelihood of transforming media to one\\u2019s liking and thereby causing psychological manipulation of society as a whole.\",\n  ],\n\n  // Paragraph: \"Since the very beginning of multimedia, ...\"\n  [\n    \"people have always tried to fake genuine information for fulfilling deceptive motives.\",\n    \"people have always tried to genuine fake information to fulfilling deceptive motives.\",\n  ],\n  [\n    \"a high potential to deceive general audience.\",\n    \"a high potential to deceive the general audience.\",\n  ],\n  [\n    \"feel of authenticity, is what makes\",\n    \"feel of authenticity is what makes\",\n  ],\n\n  // Paragraph: \"The biggest dilemma that arises due to this is, ...\"\n  [\n    \"with faces of our know and trusted people, are they feeding us with hoax?\",\n    \"with faces of our knowledgeable and trusted people, are they feeding us with the hoax?\",\n  ],\n\n  // Paragraph: \"Like most Machine Learning algorithms, ...\"\n  [\n    \"usually training into an Autoencoders or a Generative\",\n    \"usually training into Autoencoders or a Generative\",\n  ],\n\n  // Paragraph: \"On the other hand, researchers in big tech firms ...\"\n  [\n    \"Kaggle and Drivendata are working on their own machine learning perspectives to device a counteracting\",\n    \"Kaggle and Driven data are working on their machine learning perspectives to devise a counteracting\",\n  ],\n\n  // Paragraph: \"Many other technologies are on the verge ...\"\n  [\n    \"videos might be so perfect that its physically impossible to distinguish.\",\n    \"videos might be so perfect that it's physically impossible to distinguish.\",\n  ],\n\n  // Paragraph: \"What we need to do is to trust the source ...\"\n  [\n    \"One of the major causes of spread of misleading information, is the fact that\",\n    \"One of the major causes of the spread of misleading information is the fact that\",\n  ],\n  [\n    \"keep our eyes and ears open for anyone trying to invasively shape our believes.\",\n    \"keep our eyes and ears open for anyone trying to invasively shape our beliefs.\",\n  ],\n];\n\nfor (const [searchText, replacement] of edits) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + searchText);\n  }\n\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word / $d resolve against the already-open document (ActiveDocument).\n# We locate each target phrase with Find.Execute (search only, no\n# auto-replace) and then assign the new wording directly via Range.Text,\n# which edits the text in place without Word's Find/Replace \"smart quotes\"\n# AutoCorrect post-processing silently changing straight apostrophes.\n\n$d = $word.ActiveDocument\n\nfunction Set-RangeText($findText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap(wdFindContinue=1), Format,\n    # ReplaceWith, Replace(wdReplaceNone=0 -- locate only, don't auto-replace)\n    $found = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n    $r = $find.Parent\n    $r.Text = $newText\n}\n\n$curly = [char]0x2019\n\n# Paragraph: \"Everything we hear or see, ...\"\nSet-RangeText \"school, or that burst\" \"school or that burst\"\nSet-RangeText \"significantly, and has leveraged our dependence of technology, completely ignoring the down-sides\" \"significantly and has leveraged our dependence on technology, completely ignoring the downsides\"\nSet-RangeText \"All of this, increases the likelihood of transforming media into ones liking and thereby causing psychological manipulation of the society as a whole.\" (\"All of this increases the likelihood of transforming media to one\" + $curly + \"s liking and thereby causing psychological manipulation of society as a whole.\")\n\n# Paragraph: \"Since the very beginning of multimedia, ...\"\nSet-RangeText \"people have always tried to fake genuine information for fulfilling deceptive motives.\" \"people have always tried to genuine fake information to fulfilling deceptive motives.\"\nSet-RangeText \"a high potential to deceive general audience.\" \"a high potential to deceive the general audience.\"\nSet-RangeText \"feel of authenticity, is what makes\" \"feel of authenticity is what makes\"\n\n# Paragraph: \"The biggest dilemma that arises due to this is, ...\"\nSet-RangeText \"with faces of our know and trusted people, are they feeding us with hoax?\" \"with faces of our knowledgeable and trusted people, are they feeding us with the hoax?\"\n\n# Paragraph: \"Like most Machine Learning algorithms, ...\"\nSet-RangeText \"usually training into an Autoencoders or a Generative\" \"usually training into Autoencoders or a Generative\"\n\n# Paragraph: \"On the other hand, researchers in big tech firms ...\"\nSet-RangeText \"Kaggle and Drivendata are working on their own machine learning perspectives to device a counteracting\" \"Kaggle and Driven data are working on their machine learning perspectives to devise a counteracting\"\n\n# Paragraph: \"Many other technologies are on the verge ...\"\nSet-RangeText \"videos might be so perfect that its physically impossible to distinguish.\" \"videos might be so perfect that it's physically impossible to distinguish.\"\n\n# Paragraph: \"What we need to do is to trust the source ...\"\nSet-RangeText \"One of the major causes of spread of misleading information, is the fact that\" \"One of the major causes of the spread of misleading information is the fact that\"\nSet-RangeText \"keep our eyes and ears open for anyone trying to invasively shape our believes.\" \"keep our eyes and ears open for anyone trying to invasively shape our beliefs.\"\n"}
